# Munkfors "Avverkningsanmälningar" overview update.
#
# A new felling-notice record (A 44111-2023) was added as a new row 2,
# pushing every existing data row (formerly rows 2-84) down by one (now
# rows 3-85). In addition every pre-existing record's "Förändrad" date
# (column C) was refreshed to the new run date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row and give it the same row height as its neighbours ---
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).RowHeight = 15

# --- Populate the new row 2 with the new record ---
$ws.Range("A2").Value = "A 44111-2023"

$ws.Range("B2").Value = 45188
$ws.Range("B2").NumberFormat = "YYYY-MM-DD"

$ws.Range("C2").Value = 45190
$ws.Range("C2").NumberFormat = "YYYY-MM-DD"

$ws.Range("D2").Value = "VÄRMLANDS LÄN"
$ws.Range("E2").Value = "MUNKFORS"
# F2 (Markägare) intentionally left blank for this record.

$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 1

$ws.Range("R2").Value = "Bågsäv"
$ws.Range("R2").WrapText = $true

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MUNKFORS/artfynd/A 44111-2023.xlsx", "A 44111-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MUNKFORS/kartor/A 44111-2023.png", "A 44111-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MUNKFORS/klagomål/A 44111-2023.docx", "A 44111-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MUNKFORS/klagomålsmail/A 44111-2023.docx", "A 44111-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MUNKFORS/tillsyn/A 44111-2023.docx", "A 44111-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MUNKFORS/tillsynsmail/A 44111-2023.docx", "A 44111-2023")'

# --- Refresh the "Förändrad" date on every pre-existing record (now rows 3-85) ---
$ws.Range("C3:C85").Value = 45190
